$d = $word.ActiveDocument

# Update the date line at the top of the document
$d.Content.Find.Execute("2026-02-11 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2026-02-12 Thursday", 2)

# Update each multiplication problem cell in the table (targeted by row/col to avoid
# ambiguity, since some new values coincide with other cells' old values)
$t = $d.Tables.Item(1)
$t.Cell(1,1).Range.Find.Execute("615×8=", $true, $false, $false, $false, $false, $true, 1, $false, "809×7=", 2)
$t.Cell(1,2).Range.Find.Execute("991×2=", $true, $false, $false, $false, $false, $true, 1, $false, "516×2=", 2)
$t.Cell(1,3).Range.Find.Execute("405×2=", $true, $false, $false, $false, $false, $true, 1, $false, "174×7=", 2)
$t.Cell(1,4).Range.Find.Execute("202×6=", $true, $false, $false, $false, $false, $true, 1, $false, "419×3=", 2)
$t.Cell(1,5).Range.Find.Execute("710×6=", $true, $false, $false, $false, $false, $true, 1, $false, "799×9=", 2)
$t.Cell(5,1).Range.Find.Execute("454×2=", $true, $false, $false, $false, $false, $true, 1, $false, "548×7=", 2)
$t.Cell(5,2).Range.Find.Execute("432×4=", $true, $false, $false, $false, $false, $true, 1, $false, "188×7=", 2)
$t.Cell(5,3).Range.Find.Execute("439×8=", $true, $false, $false, $false, $false, $true, 1, $false, "352×3=", 2)
$t.Cell(5,4).Range.Find.Execute("444×3=", $true, $false, $false, $false, $false, $true, 1, $false, "543×5=", 2)
$t.Cell(5,5).Range.Find.Execute("125×9=", $true, $false, $false, $false, $false, $true, 1, $false, "222×6=", 2)
$t.Cell(10,1).Range.Find.Execute("609×6=", $true, $false, $false, $false, $false, $true, 1, $false, "422×6=", 2)
$t.Cell(10,2).Range.Find.Execute("371×6=", $true, $false, $false, $false, $false, $true, 1, $false, "736×5=", 2)
$t.Cell(10,3).Range.Find.Execute("742×4=", $true, $false, $false, $false, $false, $true, 1, $false, "691×9=", 2)
$t.Cell(10,4).Range.Find.Execute("809×5=", $true, $false, $false, $false, $false, $true, 1, $false, "439×8=", 2)
$t.Cell(10,5).Range.Find.Execute("269×4=", $true, $false, $false, $false, $false, $true, 1, $false, "161×9=", 2)
$t.Cell(15,1).Range.Find.Execute("614×9=", $true, $false, $false, $false, $false, $true, 1, $false, "850×2=", 2)
$t.Cell(15,2).Range.Find.Execute("597×5=", $true, $false, $false, $false, $false, $true, 1, $false, "297×9=", 2)
$t.Cell(15,3).Range.Find.Execute("552×4=", $true, $false, $false, $false, $false, $true, 1, $false, "929×9=", 2)
$t.Cell(15,4).Range.Find.Execute("158×8=", $true, $false, $false, $false, $false, $true, 1, $false, "147×3=", 2)
$t.Cell(15,5).Range.Find.Execute("218×2=", $true, $false, $false, $false, $false, $true, 1, $false, "914×6=", 2)
$t.Cell(20,1).Range.Find.Execute("510×2=", $true, $false, $false, $false, $false, $true, 1, $false, "596×8=", 2)
$t.Cell(20,2).Range.Find.Execute("455×6=", $true, $false, $false, $false, $false, $true, 1, $false, "139×5=", 2)
$t.Cell(20,3).Range.Find.Execute("207×9=", $true, $false, $false, $false, $false, $true, 1, $false, "750×3=", 2)
$t.Cell(20,4).Range.Find.Execute("172×5=", $true, $false, $false, $false, $false, $true, 1, $false, "180×7=", 2)
$t.Cell(20,5).Range.Find.Execute("246×9=", $true, $false, $false, $false, $false, $true, 1, $false, "520×6=", 2)
